$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Great work, continue with excellence!"
$ws.Range("C3").Value = "My name is Daniel."
$ws.Range("C4").Value = "Hello, this is a test using the OpenAI API to correct text."
$ws.Range("C5").Value = "I believe in you."
$ws.Range("C6").Value = "You have a bright future ahead. Keep moving forward."
$ws.Range("C7").Value = "You are capable of doing incredible things."
$ws.Range("C8").Value = "Your performance is improving every day."
$ws.Range("C9").Value = "You are intelligent and I am proud of it."
$ws.Range("C10").Value = "Your effort is worth it."
$ws.Range("C11").Value = "You are an inspiration."
